$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "-"
$ws.Range("B3").Value = "[-, -, -, 'MEC-3B-Fresagem']"
$ws.Range("E3").Value = "-"
$ws.Range("B4").Value = "[-, -, -, 'MEC-3B-Fresagem']"
$ws.Range("E4").Value = "-"
$ws.Range("B6").Value = "[-, -, -, 'MEC-3B-Fresagem']"
$ws.Range("B7").Value = "[-, -, -, 'MEC-3B-Fresagem']"
$ws.Range("E8").Value = "-"
$ws.Range("F10").Value = "-"
$ws.Range("E11").Value = "['MEC-3A-Fresagem', -, -, -]"
$ws.Range("F11").Value = "[-, -, 'MEC-2A-Tornearia', -]"
$ws.Range("E12").Value = "['MEC-3A-Fresagem', -, -, -]"
$ws.Range("F12").Value = "[-, -, 'MEC-2A-Tornearia', -]"
$ws.Range("E14").Value = "['MEC-3A-Fresagem', -, -, -]"
$ws.Range("F14").Value = "[-, -, 'MEC-2A-Tornearia', -]"
$ws.Range("E15").Value = "['MEC-3A-Fresagem', -, -, -]"
$ws.Range("F15").Value = "[-, -, 'MEC-2A-Tornearia', -]"
$ws.Range("F16").Value = "-"
$ws.Range("B18").Value = "[-, -, -, 'MEC-2NA-Fresagem']"
$ws.Range("C18").Value = "[-, -, -, 'MEC-2NA-CAD/CAM']"
$ws.Range("D18").Value = "[-, -, -, 'MEC-2NA-CAD/CAM']"
$ws.Range("E18").Value = "[-, -, -, 'MEC-2NB-Fresagem']"
$ws.Range("B19").Value = "[-, -, 'MEC-2NA-CAD/CAM', 'MEC-2NA-Fresagem']"
$ws.Range("C19").Value = "-"
$ws.Range("D19").Value = "-"
$ws.Range("E19").Value = "[-, -, -, 'MEC-2NB-Fresagem']"
$ws.Range("B20").Value = "[-, -, 'MEC-2NA-CAD/CAM', 'MEC-2NA-Fresagem']"
$ws.Range("C20").Value = "-"
$ws.Range("D20").Value = "-"
$ws.Range("E20").Value = "[-, -, -, 'MEC-2NB-Fresagem']"
$ws.Range("B21").Value = "[-, -, -, 'MEC-2NA-Fresagem']"
$ws.Range("C21").Value = "-"
$ws.Range("D21").Value = "-"
$ws.Range("E21").Value = "[-, -, -, 'MEC-2NB-Fresagem']"
